$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - reorder "Recorded By" list
$ws.Range("G2").Value = "servinaz@med.asu.edu.eg, System, gehanadel@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 3 - reorder "Recorded By" list
$ws.Range("G3").Value = "hend_mahmoud@med.asu.edu.eg, System, asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 4 - reorder "Recorded By" list
$ws.Range("G4").Value = "hend_mahmoud@med.asu.edu.eg, servinaz@med.asu.edu.eg, asmaa.reda@med.asu.edu.eg, gehanadel@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, majorelle.magdy@med.asu.edu.eg"

# Row 5 - reorder "Recorded By" list
$ws.Range("G5").Value = "asmaa.reda@med.asu.edu.eg, eman.tantawi@med.asu.edu.eg, Amira.Sobhy@med.asu.edu.eg, Veronia.rafat@med.asu.edu.eg"

# Row 6 - reorder "Recorded By" list, Missing->Recorded count bump
$ws.Range("G6").Value = "alshimaa.atef@med.asu.edu.egm, majorelle.magdy@med.asu.edu.eg, manar.montaser@med.asu.edu.eg, Mohammedeltanany@med.asu.edu.eg"
$ws.Range("L6").Value = 27

# Row 7 - reorder "Recorded By" list, Missing Sessions count drop
$ws.Range("G7").Value = "Amera.a.saad@med.asu.edu.eg, Fatmaelhady@med.asu.edu.eg, AbeerRagheb@med.asu.edu.eg, NadaMohamed@med.asu.edu.eg, lamiaa.ossama@med.asu.edu.eg, menna-alah.mohamed@asu.edu.eg, Kerelos.zareef@med.asu.edu.eg"
$ws.Range("L7").Value = 2

# Row 9 - Coverage % stat update (force literal text so it isn't
# auto-converted to a percentage number; then restore the original cell
# formatting that NumberFormat="@" would otherwise overwrite)
$ws.Range("L9").NumberFormat = "@"
$ws.Range("L9").Value = "93.1%"
$ws.Range("K9").Copy()
$ws.Range("L9").PasteSpecial(-4122)

# Row 10 - Average Attendance % stat update
$ws.Range("L10").NumberFormat = "@"
$ws.Range("L10").Value = "26.9%"
$ws.Range("K10").Copy()
$ws.Range("L10").PasteSpecial(-4122)

# Row 12 - reorder "Recorded By" list
$ws.Range("G12").Value = "Madeha.Saeed@med.asu.edu.eg, yassmina.fattoh@med.asu.edu.eg, dina.adel@med.asu.edu.eg, amira.m.ibrahim@med.asu.edu.eg, Eman.m.abosakaya@med.asu.edu.eg, Marina.youhana@med.asu.edu.eg"

# Row 15 - summary table stats (Recorded/Missing counts + percentages)
$ws.Range("O15").Value = 27
$ws.Range("P15").Value = 2

$ws.Range("R15").NumberFormat = "@"
$ws.Range("R15").Value = "93.1%"
$ws.Range("Q15").Copy()
$ws.Range("R15").PasteSpecial(-4122)

$ws.Range("S15").NumberFormat = "@"
$ws.Range("S15").Value = "26.9%"
$ws.Range("Q15").Copy()
$ws.Range("S15").PasteSpecial(-4122)

# Row 17 - reorder "Recorded By" list
$ws.Range("G17").Value = "mohamed.saleem@med.asu.edu.eg, esraa.sami@med.asu.edu.eg"

# Row 19 - session now recorded: copy "Recorded" formatting from row 18, then set values
$ws.Range("A18:I18").Copy()
$ws.Range("A19:I19").PasteSpecial(-4122)
$ws.Range("G19").Value = "Rania.a.youssef@med.asu.edu.eg"
$ws.Range("H19").Value = "18/251"
$ws.Range("I19").Value = "Recorded"

# Row 20 - reorder "Recorded By" list
$ws.Range("G20").Value = "mohamed.saleem@med.asu.edu.eg, mariam.youssif.std@med.asu.edu.eg"

# Row 25 - reorder "Recorded By" list
$ws.Range("G25").Value = "Noran.Mahmoud@med.asu.edu.eg, menna-allah.gamil@med.asu.edu.eg"

# Row 27 - reorder "Recorded By" list
$ws.Range("G27").Value = "nourhan.mostafa@med.asu.edu.eg, hana.amr@med.asu.edu.eg"

# Row 28 - reorder "Recorded By" list
$ws.Range("G28").Value = "Aya_hamed@med.asu.edu.eg, maryam.ashraf@med.asu.edu.eg"

# Row 30 - reorder "Recorded By" list
$ws.Range("G30").Value = "aya.hanafy@med.asu.edu.eg, wafaa.ebida@med.asu.edu.eg, yassmen.ahmed@med.asu.edu.eg, shorokmohamed@med.asu.edu.eg"
